# Update the stock symbol table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for rows 2-28 (columns A-F).
# A = index number, B..F = NSE ticker strings ("" means blank cell)
$data = @(
    @(0,  "NSE:ALPHAETF",   "NSE:ANGELONE",  "NSE:COLPAL", "", "NSE:CUB"),
    @(1,  "NSE:APOLSINHOT", "NSE:BHAGCHEM",  "NSE:OFSS",   "", "NSE:OFSS"),
    @(2,  "NSE:ARIHANTSUP", "NSE:BSE",        "", "", ""),
    @(3,  "NSE:ARVIND",     "NSE:CHOLAHLDNG", "", "", ""),
    @(4,  "NSE:BAJAJELEC",  "NSE:CROMPTON",   "", "", ""),
    @(5,  "NSE:BAJAJHLDNG", "NSE:DBSTOCKBRO", "", "", ""),
    @(6,  "NSE:CARTRADE",   "NSE:DECCANCE",   "", "", ""),
    @(7,  "NSE:CCHHL",      "NSE:EDELWEISS",  "", "", ""),
    @(8,  "NSE:CONCORDBIO", "NSE:ESCORTS",    "", "", ""),
    @(9,  "NSE:CREST",      "NSE:HLEGLAS",    "", "", ""),
    @(10, "NSE:GRAPHITE",   "NSE:IKIO",       "", "", ""),
    @(11, "NSE:GREAVESCOT", "NSE:INFIBEAM",   "", "", ""),
    @(12, "NSE:GRMOVER",    "NSE:NESCO",      "", "", ""),
    @(13, "NSE:GTLINFRA",   "NSE:PGEL",       "", "", ""),
    @(14, "NSE:GUFICBIO",   "NSE:PNBGILTS",   "", "", ""),
    @(15, "NSE:ISGEC",      "NSE:RBL",        "", "", ""),
    @(16, "NSE:KIOCL",      "", "", "", ""),
    @(17, "NSE:LICMFGOLD",  "", "", "", ""),
    @(18, "NSE:LINDEINDIA", "", "", "", ""),
    @(19, "NSE:MANGLMCEM",  "", "", "", ""),
    @(20, "NSE:NBCC",       "", "", "", ""),
    @(21, "NSE:NCLIND",     "", "", "", ""),
    @(22, "NSE:NOCIL",      "", "", "", ""),
    @(23, "NSE:PITTIENG",   "", "", "", ""),
    @(24, "NSE:PLASTIBLEN", "", "", "", ""),
    @(25, "NSE:QGOLDHALF",  "", "", "", ""),
    @(26, "NSE:RATEGAIN",   "", "", "", "")
)

# First, delete the rows that are no longer needed (rows 29-37), shifting
# everything below up. Doing this first keeps row numbers 2-28 stable while
# we overwrite their contents below.
$ws.Range("A29:F37").EntireRow.Delete() | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($rowNum, 1).Value = $row[0]
    $ws.Cells.Item($rowNum, 2).Value = $row[1]
    $ws.Cells.Item($rowNum, 3).Value = $row[2]
    $ws.Cells.Item($rowNum, 4).Value = $row[3]
    $ws.Cells.Item($rowNum, 5).Value = $row[4]
    $ws.Cells.Item($rowNum, 6).Value = $row[5]
}
